# Saldo.xlsx update: remove three stale "Export" rows and add a new one.
#
# Rows removed (matched by account number in column A):
#   000834301 | MARCUS   | 104229.07
#   004752519 | MARCUS   |  46399.01
#   004388077 | WLADMIR  |  39673.18
#
# Row added (inserted right after account 004454365 / RAFAEL / 13566.65,
# i.e. right before account 004487140 / VALMIR):
#   005587298 | JULIO    |  10000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from the bottom up so earlier row numbers stay valid while we work.
$ws.Rows.Item(9).Delete()   # 004388077 | WLADMIR  | 39673.18
$ws.Rows.Item(7).Delete()   # 004752519 | MARCUS   | 46399.01
$ws.Rows.Item(5).Delete()   # 000834301 | MARCUS   | 104229.07

# After the three deletions above, row 9 is 004454365/RAFAEL/13566.65 and
# row 10 is 004487140/VALMIR/6612. Insert the new row between them.
$ws.Rows.Item(10).Insert()

# Column A holds zero-padded account numbers stored as text; format the new
# cell as Text first so the leading zeros in "005587298" are preserved.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "005587298"
$ws.Range("B10").Value = "JULIO"
$ws.Range("C10").Value = 10000
